$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A is formatted as Text before any values are entered (so both the
# column default style and the individual cells pick up numFmtId 49).
$ws.Columns.Item(1).NumberFormat = "@"
$ws.Columns.Item(1).ColumnWidth = 17.333333333333336

# --- Header row (row 1) -----------------------------------------------
$ws.Range("A1").Value = "TGL.TRANSAKSI"
$ws.Range("B1").Value = "JNS.TRANSAKSI"
$ws.Range("C1").Value = "NAMA OBAT/ ALKES"
$ws.Range("D1").Value = "DEPO FARMASI"
$ws.Range("F1").Value = "DEBET"
$ws.Range("G1").Value = "KREDIT"
$ws.Range("H1").Value = "BIAYA"

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = "2019-09-06 05:45:57"
$ws.Range("C2").Value = "AMOXYCILIN 500 MG"
$ws.Range("B2").Value = "PEMBELIAN"
$ws.Range("D2").Value = "GUDANG OBAT"
$ws.Range("F2").Value = 70
$ws.Range("H2").Value = 17149671000

# A2 additionally gets an explicit black font colour (creates a second,
# text-formatted font/style combination).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Font.Color = 0

# --- Column E (KATEGORI / ASKES) added last --------------------------
$ws.Range("E1").Value = "KATEGORI"
$ws.Range("E2").Value = "ASKES"

# --- Column widths for the rest of the used range ----------------------
# (Values are chosen so the engine's column-width quantisation lands as
# close as possible to the widths Excel's own AutoFit produced: 14.7109375,
# 18.85546875, 14.42578125, 13.5703125, 11.140625, 12, 12.140625.)
$ws.Columns.Item(2).ColumnWidth = 13.833333333333332
$ws.Columns.Item(3).ColumnWidth = 18.0
$ws.Columns.Item(4).ColumnWidth = 13.666666666666668
$ws.Columns.Item(5).ColumnWidth = 12.666666666666668
$ws.Columns.Item(7).ColumnWidth = 10.333333333333332
$ws.Columns.Item(8).ColumnWidth = 11.166666666666668
$ws.Columns.Item(9).ColumnWidth = 11.333333333333332

$ws.Range("E2").Select()
